$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.03069236259814418
$ws.Range("D3").Value = 0.8431372549019608
$ws.Range("E3").Value = 0.05922865013774105

$ws.Range("C4").Value = 0.04298356510745891
$ws.Range("D4").Value = 0.6666666666666666
$ws.Range("E4").Value = 0.08076009501187648

$ws.Range("C5").Value = 0.04627539503386004
$ws.Range("D5").Value = 0.803921568627451
$ws.Range("E5").Value = 0.08751334044823907

$ws.Range("C6").Value = 0.06551724137931035
$ws.Range("D6").Value = 0.3725490196078431
$ws.Range("E6").Value = 0.1114369501466276

$ws.Range("C7").Value = 0.07888631090487239
$ws.Range("D7").Value = 0.6666666666666666
$ws.Range("E7").Value = 0.1410788381742738

$ws.Range("C8").Value = 0.0583941605839416
$ws.Range("D8").Value = 0.4705882352941176
$ws.Range("E8").Value = 0.1038961038961039
